$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.184.22"
$ws.Range("E2").Value = "  +2.95%  "

# Row 3
$ws.Range("D3").Value = "3.483.99"
$ws.Range("E3").Value = "  +2.58%  "

# Row 5
$ws.Range("D5").Value = "'579.81"
$ws.Range("E5").Value = "  +2.11%  "

# Row 6
$ws.Range("D6").Value = "'162.82"
$ws.Range("E6").Value = "  +4.65%  "

# Row 7
$ws.Range("D7").Value = "'0.615"
$ws.Range("E7").Value = "  +12.68%  "

# Row 8
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").Value = "3.487.44"
$ws.Range("E9").Value = "  +2.68%  "

# Row 10
$ws.Range("D10").Value = "'7.26"
$ws.Range("E10").Value = "  -1.92%  "

# Row 11
$ws.Range("D11").Value = "'0.125"
$ws.Range("E11").Value = "  +3.00%  "

# Row 12
$ws.Range("D12").Value = "'0.447"
$ws.Range("E12").Value = "  +3.73%  "

# Row 13
$ws.Range("D13").Value = "4.089.79"
$ws.Range("E13").Value = "  +2.65%  "

# Row 14
$ws.Range("D14").Value = "'0.135"
$ws.Range("E14").Value = "  +0.55%  "

# Row 15
$ws.Range("D15").Value = "'0.0000193"
$ws.Range("E15").Value = "  +0.09%  "

# Row 16
$ws.Range("D16").Value = "'28.70"
$ws.Range("E16").Value = "  +5.40%  "

# Row 17
$ws.Range("D17").Value = "65.242.69"
$ws.Range("E17").Value = "  +2.95%  "

# Row 18
$ws.Range("D18").Value = "3.496.19"
$ws.Range("E18").Value = "  +2.50%  "

# Row 19
$ws.Range("D19").Value = "'6.47"
$ws.Range("E19").Value = "  +3.68%  "

# Row 20
$ws.Range("D20").Value = "'14.38"
$ws.Range("E20").Value = "  +2.24%  "

# Row 21
$ws.Range("D21").Value = "'383.41"
$ws.Range("E21").Value = "  +0.95%  "

# Row 22
$ws.Range("D22").Value = "'8.20"
$ws.Range("E22").Value = "  +2.04%  "

# Row 23
$ws.Range("D23").Value = "'0.554"
$ws.Range("E23").Value = "  +4.71%  "

# Row 24
$ws.Range("D24").Value = "'72.66"
$ws.Range("E24").Value = "  +1.62%  "

# Row 25
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.11%  "

# Row 26
$ws.Range("D26").Value = "'0.0000120"
$ws.Range("E26").Value = "  +1.07%  "

# Row 27
$ws.Range("D27").Value = "'10.03"
$ws.Range("E27").Value = "  +6.55%  "

# Row 28
$ws.Range("E28").Value = "  +0.21%  "

# Row 29
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.03%  "

# Row 30
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.54"
$ws.Range("E30").Value = "  +13.53%  "

# Row 31
$ws.Range("D31").Value = "'6.16"
$ws.Range("E31").Value = "  +1.45%  "

# Row 32
$ws.Range("D32").Value = "'2.06"
$ws.Range("E32").Value = "  +2.89%  "

# Row 33
$ws.Range("D33").Value = "'23.71"
$ws.Range("E33").Value = "  +2.37%  "

# Row 34
$ws.Range("D34").Value = "'7.19"
$ws.Range("E34").Value = "  +5.75%  "

# Row 35
$ws.Range("D35").Value = "'1.63"
$ws.Range("E35").Value = "  +12.44%  "

# Row 36
$ws.Range("D36").Value = "'162.21"
$ws.Range("E36").Value = "  +1.47%  "

# Row 37
$ws.Range("E37").Value = "  +5.25%  "

# Row 38
$ws.Range("D38").Value = "'0.0780"
$ws.Range("E38").Value = "  +4.00%  "

# Row 39
$ws.Range("D39").Value = "2.999.75"
$ws.Range("E39").Value = "  +1.74%  "

# Row 40
$ws.Range("D40").Value = "'6.83"
$ws.Range("E40").Value = "  +7.49%  "

# Row 41
$ws.Range("D41").Value = "'26.83"
$ws.Range("E41").Value = "  -0.57%  "

# Row 42
$ws.Range("D42").Value = "'4.58"
$ws.Range("E42").Value = "  +5.78%  "

# Row 43
$ws.Range("E43").Value = "  +2.03%  "

# Row 44
$ws.Range("D44").Value = "'42.91"
$ws.Range("E44").Value = "  +2.78%  "

# Row 45
$ws.Range("D45").Value = "'0.782"
$ws.Range("E45").Value = "  +2.64%  "

# Row 46
$ws.Range("D46").Value = "'25.88"
$ws.Range("E46").Value = "  +11.22%  "

# Row 47
$ws.Range("D47").Value = "'1.11"
$ws.Range("E47").Value = "  +3.75%  "

# Row 48
$ws.Range("D48").Value = "'319.74"
$ws.Range("E48").Value = "  +9.69%  "

# Row 49
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'6.73"
$ws.Range("E49").Value = "  +5.97%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.110"
$ws.Range("E50").Value = "  +6.54%  "

# Row 51
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "'0.881"
$ws.Range("E51").Value = "  +5.07%  "
